# Split the single run "demo" into four runs: "D" | "emo" | " " | "4"
# (all sharing the original <w:lang w:val="en-US"/> run formatting), so the
# paragraph ends up reading "Demo 4" but spread across separate <w:r> nodes
# exactly like the target OOXML diff.

$d = $word.ActiveDocument

# Locate the exact text we need to replace ("demo") without disturbing
# anything else in the document.
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Execute("demo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if (-not $find.Find.Found) {
    throw "Could not find target text 'demo' in document."
}

$target = $d.Range($find.Start, $find.End)

# Build a minimal WordprocessingML package fragment describing the four
# runs we want in place of the single "demo" run. Using InsertXML (instead
# of Range.Text / InsertAfter / InsertBefore) preserves the run boundaries
# instead of Word's usual same-formatting run coalescing.
$openXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>D</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>emo</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>4</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($openXml)
